# Update the F-column (time_taken) timestamps on the "data" sheet to reflect
# the refined/re-run panel query times.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:22:38.853154"
$ws1.Range("F3").Value = "2021-10-05 14:22:38.853162"
$ws1.Range("F4").Value = "2021-10-05 14:22:38.853165"
$ws1.Range("F5").Value = "2021-10-05 14:22:38.853168"
$ws1.Range("F6").Value = "2021-10-05 14:22:38.853171"
$ws1.Range("F7").Value = "2021-10-05 14:22:38.853173"
$ws1.Range("F8").Value = "2021-10-05 14:22:38.853176"
$ws1.Range("F9").Value = "2021-10-05 14:22:38.853178"
$ws1.Range("F10").Value = "2021-10-05 14:22:38.853181"
$ws1.Range("F11").Value = "2021-10-05 14:22:38.853184"
$ws1.Range("F12").Value = "2021-10-05 14:22:38.853186"
$ws1.Range("F13").Value = "2021-10-05 14:22:38.853189"
$ws1.Range("F14").Value = "2021-10-05 14:22:38.853191"
$ws1.Range("F15").Value = "2021-10-05 14:22:38.853194"
$ws1.Range("F16").Value = "2021-10-05 14:22:38.853196"
$ws1.Range("F17").Value = "2021-10-05 14:22:38.853199"
$ws1.Range("F18").Value = "2021-10-05 14:22:38.853202"
$ws1.Range("F19").Value = "2021-10-05 14:22:38.853204"
$ws1.Range("F20").Value = "2021-10-05 14:22:38.853207"
$ws1.Range("F21").Value = "2021-10-05 14:22:38.853209"
$ws1.Range("F22").Value = "2021-10-05 14:22:38.853212"
$ws1.Range("F23").Value = "2021-10-05 14:22:38.853214"
$ws1.Range("F24").Value = "2021-10-05 14:22:38.853216"
$ws1.Range("F25").Value = "2021-10-05 14:22:38.853219"
$ws1.Range("F26").Value = "2021-10-05 14:22:38.853222"
$ws1.Range("F27").Value = "2021-10-05 14:22:38.853224"
$ws1.Range("F28").Value = "2021-10-05 14:22:38.853227"
$ws1.Range("F29").Value = "2021-10-05 14:22:38.853229"
$ws1.Range("F30").Value = "2021-10-05 14:22:38.853232"
$ws1.Range("F31").Value = "2021-10-05 14:22:38.853234"
$ws1.Range("F32").Value = "2021-10-05 14:22:38.853237"
$ws1.Range("F33").Value = "2021-10-05 14:22:38.853239"
$ws1.Range("F34").Value = "2021-10-05 14:22:38.853242"
$ws1.Range("F35").Value = "2021-10-05 14:22:38.853245"
$ws1.Range("F36").Value = "2021-10-05 14:22:38.853247"
$ws1.Range("F37").Value = "2021-10-05 14:22:38.853250"
$ws1.Range("F38").Value = "2021-10-05 14:22:38.853252"
$ws1.Range("F39").Value = "2021-10-05 14:22:38.853255"
$ws1.Range("F40").Value = "2021-10-05 14:22:38.853257"
$ws1.Range("F41").Value = "2021-10-05 14:22:38.853260"
$ws1.Range("F42").Value = "2021-10-05 14:22:38.853263"

# Add a new "metadata" worksheet right after the "data" sheet, carrying the
# panel query metadata that used to live alongside the gene data.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "metadata"

# Copy the header/first-data-row formatting (bold, bordered, centered style)
# from the "data" sheet so the new sheet's header row + index column match
# the existing look-and-feel.
$ws1.Range("A1:F2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)

# Header row.
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Severe early-onset obesity"
$newSheet.Range("C2").Value = 130

# data_version must stay text ("2.43"), not be coerced to a number.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.43"

$newSheet.Range("E2").Value = "2021-09-28T15:05:45.249854Z"
$newSheet.Range("F2").Value = "2021-10-05 14:22:38.849463"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/130/?format=json"
